# Thesis data workbook update:
#  - Insert a new "Notes" worksheet before the existing "Data" sheet,
#    describing the fields used on the Data sheet.
#  - Clear out the actual student records from the "Data" sheet (rows 2-5),
#    leaving only the header row and the pre-existing blank/style rows.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")

# --- Create the new "Notes" sheet immediately before "Data" ---
$notesSheet = $wb.Worksheets.Add($dataSheet)
$notesSheet.Name = "Notes"

$notesSheet.Range("A1").Value = "Notes"

$notesSheet.Range("A3").Value = "Fields"

$notesSheet.Range("A4").Value = "Student"
$notesSheet.Range("B4").Value = "Student Name, any format works but last, first is probably best"

$notesSheet.Range("A5").Value = "Start Date"
$notesSheet.Range("B5").Value = "Date started in program (not used)"

$notesSheet.Range("A6").Value = "Year"
$notesSheet.Range("B6").Value = "Graduation year"

$notesSheet.Range("A7").Value = "Degree"
$notesSheet.Range("B7").Value = "Typically M.S. or Ph.D."

$notesSheet.Range("A8").Value = "Advisor"
$notesSheet.Range("B8").Value = "Your last name (not used)"

$notesSheet.Range("A9").Value = "Title"
$notesSheet.Range("B9").Value = 'Title of thesis.  If coadvised, put the co-advisors names here in parentheses i.e. "Title (co-advised w/J. Doe)"'

$notesSheet.Range("A10").Value = "Comments"
$notesSheet.Range("B10").Value = "Not used."

# --- Clear the old student records from the "Data" sheet ---
# Re-fetch the "Data" sheet reference by name: inserting the new sheet
# shifted its position, and a handle captured beforehand would now
# resolve to the new "Notes" sheet instead.
$dataSheet = $wb.Worksheets.Item("Data")
$dataSheet.Range("A2:F5").ClearContents()

# --- Update selections / active sheet to match the saved state ---
$dataSheet.Range("A1:G1").Select()
$notesSheet.Activate()
$notesSheet.Range("B18").Select()
